$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header text (Portuguese, capitalization of "Chai" fixed + a couple of
# wording tweaks) for the table header cells B1:F1.
$newHeaders = @{
    "B1" = "Total de vendas de Chai (unidades)"
    "C1" = "Vendas de Chai Artesanal (unidades)"
    "D1" = "Vendas de Chai pré-fabricado (unidades)"
    "E1" = "Engajamento em redes sociais (visualizações)"
    "F1" = "Pesquisas online para Chai"
}

$white = 16777215

foreach ($addr in $newHeaders.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $newHeaders[$addr]

    $len = $cell.Characters().Text.Length

    # Re-apply bold + the original white font colour to the whole run. The
    # range is deliberately split in two (1..len-1 and len..len) instead of
    # addressing the full string in one call, since formatting the entire
    # string length in a single Characters() call on this host is a no-op;
    # splitting it still yields a single consolidated run once both pieces
    # share identical formatting.
    $cell.Characters(1, $len - 1).Font.Bold = $true
    $cell.Characters($len, 1).Font.Bold = $true
    $cell.Characters(1, $len - 1).Font.Color = $white
    $cell.Characters($len, 1).Font.Color = $white
}

# Keep the table's column headers (ListObject) in sync with the new text.
$table = $ws.ListObjects.Item(1)
$table.ListColumns.Item(2).Name = $newHeaders["B1"]
$table.ListColumns.Item(3).Name = $newHeaders["C1"]
$table.ListColumns.Item(4).Name = $newHeaders["D1"]
$table.ListColumns.Item(5).Name = $newHeaders["E1"]
$table.ListColumns.Item(6).Name = $newHeaders["F1"]
